$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.020.75"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "1.657.34"
$ws.Range("E3").Value = "  +2.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.83"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3911"
$ws.Range("E7").Value = "  -0.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3848"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.28"
$ws.Range("E9").Value = "  +3.86%  "
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9987"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08468"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.15"
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.143"
$ws.Range("E14").Value = "  +1.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.911"
$ws.Range("E15").Value = "  +4.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001320"
$ws.Range("E16").Value = "  +3.22%  "
$ws.Range("D17").Value = "1.655.09"
$ws.Range("E17").Value = "  +2.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.73"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06977"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.84"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.933"
$ws.Range("E21").Value = "  +1.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.65"
$ws.Range("E23").Value = "  +1.72%  "
$ws.Range("D24").Value = "24.014.21"
$ws.Range("E24").Value = "  +0.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.486"
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.025"
$ws.Range("E26").Value = "  +6.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.16"
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.73"
$ws.Range("E28").Value = "  -3.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.452"
$ws.Range("E29").Value = "  +3.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "139.65"
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.893"
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.488"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").Value = "1.834.99"
$ws.Range("E33").Value = "  +1.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.043"
$ws.Range("E34").Value = "  +7.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08116"
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02988"
$ws.Range("E36").Value = "  +3.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.776"
$ws.Range("E37").Value = "  +3.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.90"
$ws.Range("E38").Value = "  +5.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2686"
$ws.Range("E39").Value = "  +0.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09167"
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7588"
$ws.Range("E41").Value = "  +1.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.50"
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.41"
$ws.Range("E44").Value = "  +2.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6966"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.464"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.086"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9991"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08304"
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.76"
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.222"
$ws.Range("E51").Value = "  +1.27%  "
